$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("savings_commercial")
$added = $wb.Styles.Add("GrayPct")
Write-Host $added
$s = $wb.Styles.Item("GrayPct")
Write-Host $s
$s.NumberFormat = "0.0%"
$s.Interior.ThemeColor = 1
$ws.Range("U27").Style = $s
